# chore: adapt column header formatting to respective input file names
#
# - Renames the "*_old" / "*_new" column headers (row 1) to "*_FV2410" /
#   "*_FV2504" respectively (the "diff" header in column K is unchanged).
# - Turns the A1:U74 range into a native Excel table ("Table1") with an
#   autofilter, matching the new header names.
# - Freezes the header row (split below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename headers: "<name>_old" -> "<name>_FV2410" (columns A-J) and
#        "<name>_new" -> "<name>_FV2504" (columns L-U). Column K ("diff")
#        stays as-is.
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2504"
}

# --- 2. Convert the used range into a table so the headers double as
#        filter buttons.
$dataRange = $ws.Range("A1:U74")
$tbl = $ws.ListObjects.Add(1, $dataRange, $true)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
